# "add excel import to brute force method" - refresh the weight/reliability
# sample data used by the brute-force reliability routine and move the
# active selection to reflect where data entry continued.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# New weight (B) / reliability (C) values for rows 7-11.
$data = @(
    @{ Row = 7;  Weight = 2; Reliability = 0.7 },
    @{ Row = 8;  Weight = 6; Reliability = 0.8 },
    @{ Row = 9;  Weight = 4; Reliability = 0.6 },
    @{ Row = 10; Weight = 2; Reliability = 0.7 },
    @{ Row = 11; Weight = 6; Reliability = 0.8 }
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 2).Value = $entry.Weight
    $ws.Cells.Item($entry.Row, 3).Value = $entry.Reliability
}

# Move the active cell selection to E8.
$ws.Range("E8").Select()
